# Update "paises.xlsx" Covid dashboard: refresh the timestamp banner and
# the statistics for the countries whose case counts changed between the
# 08:22 and 08:52 pulls. Three countries (Ucrania, Afganistan, Gabon) moved
# up in the ranking with brand new numbers, which pushes the rows below
# them (previously occupied by Panama/Catar, Cuba/Tunez/Bulgaria and
# Liechtenstein..Liberia respectively) down by one - i.e. those rows now
# show the figures that used to belong to the row above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: "Datos actualizados" banner - time bumped from 08:22 to 08:52.
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 15 de Abril de 2020 a las 08:52"

# Row 31: Rumania - stats refresh (no row reshuffle).
$ws.Cells.Item(31, 1).Value = "Rumania"
$ws.Cells.Item(31, 2).Value = 6879
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 1051
$ws.Cells.Item(31, 5).Value = 5471
$ws.Cells.Item(31, 6).Value = 241
$ws.Cells.Item(31, 7).Value = 6
$ws.Cells.Item(31, 8).Value = 357

# Row 44: Ucrania moves up here with brand-new figures.
$ws.Cells.Item(44, 1).Value = "Ucrania"
$ws.Cells.Item(44, 2).Value = 3764
$ws.Cells.Item(44, 3).Value = 392
$ws.Cells.Item(44, 4).Value = 143
$ws.Cells.Item(44, 5).Value = 3513
$ws.Cells.Item(44, 6).Value = 45
$ws.Cells.Item(44, 7).Value = 10
$ws.Cells.Item(44, 8).Value = 108

# Row 45: Panama - shifted down, carries the old row-44 figures.
$ws.Cells.Item(45, 1).Value = "Panama"
$ws.Cells.Item(45, 2).Value = 3574
$ws.Cells.Item(45, 3).Value = 0
$ws.Cells.Item(45, 4).Value = 72
$ws.Cells.Item(45, 5).Value = 3407
$ws.Cells.Item(45, 6).Value = 106
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 95

# Row 46: Catar - shifted down, carries the old row-45 figures.
$ws.Cells.Item(46, 1).Value = "Catar"
$ws.Cells.Item(46, 2).Value = 3428
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = 373
$ws.Cells.Item(46, 5).Value = 3048
$ws.Cells.Item(46, 6).Value = 37
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 7

# Row 78: Oman - stats refresh (no row reshuffle).
$ws.Cells.Item(78, 1).Value = "Oman"
$ws.Cells.Item(78, 2).Value = 910
$ws.Cells.Item(78, 3).Value = 97
$ws.Cells.Item(78, 4).Value = 131
$ws.Cells.Item(78, 5).Value = 775
$ws.Cells.Item(78, 6).Value = 3
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 4

# Row 82: Afganistan moves up here with brand-new figures.
$ws.Cells.Item(82, 1).Value = "Afganistan"
$ws.Cells.Item(82, 2).Value = 784
$ws.Cells.Item(82, 3).Value = 70
$ws.Cells.Item(82, 4).Value = 43
$ws.Cells.Item(82, 5).Value = 716
$ws.Cells.Item(82, 6).Value = 0
$ws.Cells.Item(82, 7).Value = 2
$ws.Cells.Item(82, 8).Value = 25

# Row 83: Cuba - shifted down, carries the old row-82 figures.
$ws.Cells.Item(83, 1).Value = "Cuba"
$ws.Cells.Item(83, 2).Value = 766
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 4).Value = 132
$ws.Cells.Item(83, 5).Value = 613
$ws.Cells.Item(83, 6).Value = 9
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 21

# Row 84: Tunez - shifted down, carries the old row-83 figures.
$ws.Cells.Item(84, 1).Value = "Tunez"
$ws.Cells.Item(84, 2).Value = 747
$ws.Cells.Item(84, 3).Value = 0
$ws.Cells.Item(84, 4).Value = 43
$ws.Cells.Item(84, 5).Value = 670
$ws.Cells.Item(84, 6).Value = 89
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 34

# Row 85: Bulgaria - shifted down, carries the old row-84 figures.
$ws.Cells.Item(85, 1).Value = "Bulgaria"
$ws.Cells.Item(85, 2).Value = 735
$ws.Cells.Item(85, 3).Value = 22
$ws.Cells.Item(85, 4).Value = 105
$ws.Cells.Item(85, 5).Value = 594
$ws.Cells.Item(85, 6).Value = 29
$ws.Cells.Item(85, 7).Value = 1
$ws.Cells.Item(85, 8).Value = 36

# Row 121: Islas Feroe - stats refresh (no row reshuffle).
$ws.Cells.Item(121, 1).Value = "Islas Feroe"
$ws.Cells.Item(121, 2).Value = 184
$ws.Cells.Item(121, 3).Value = 0
$ws.Cells.Item(121, 4).Value = 166
$ws.Cells.Item(121, 5).Value = 18
$ws.Cells.Item(121, 6).Value = 0
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = 0

# Row 139: Gabon moves up here with brand-new figures.
$ws.Cells.Item(139, 1).Value = "Gabon"
$ws.Cells.Item(139, 2).Value = 80
$ws.Cells.Item(139, 3).Value = 23
$ws.Cells.Item(139, 4).Value = 4
$ws.Cells.Item(139, 5).Value = 75
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 1

# Row 140: Liechtenstein - shifted down, carries the old row-139 figures.
$ws.Cells.Item(140, 1).Value = "Liechtenstein"
$ws.Cells.Item(140, 2).Value = 79
$ws.Cells.Item(140, 3).Value = 0
$ws.Cells.Item(140, 4).Value = 55
$ws.Cells.Item(140, 5).Value = 23
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 1

# Row 141: Togo - shifted down, carries the old row-140 figures.
$ws.Cells.Item(141, 1).Value = "Togo"
$ws.Cells.Item(141, 2).Value = 77
$ws.Cells.Item(141, 3).Value = 0
$ws.Cells.Item(141, 4).Value = 32
$ws.Cells.Item(141, 5).Value = 42
$ws.Cells.Item(141, 6).Value = 0
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = 3

# Row 142: Birmania - shifted down, carries the old row-141 figures.
$ws.Cells.Item(142, 1).Value = "Birmania"
$ws.Cells.Item(142, 2).Value = 74
$ws.Cells.Item(142, 3).Value = 11
$ws.Cells.Item(142, 4).Value = 2
$ws.Cells.Item(142, 5).Value = 68
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(142, 8).Value = 4

# Row 143: Congo - shifted down, carries the old row-142 figures.
$ws.Cells.Item(143, 1).Value = "Congo"
$ws.Cells.Item(143, 2).Value = 74
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 10
$ws.Cells.Item(143, 5).Value = 59
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 5

# Row 144: Barbados - shifted down, carries the old row-143 figures.
$ws.Cells.Item(144, 1).Value = "Barbados"
$ws.Cells.Item(144, 2).Value = 73
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 15
$ws.Cells.Item(144, 5).Value = 53
$ws.Cells.Item(144, 6).Value = 4
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 5

# Row 145: Somalia - shifted down, carries the old row-144 figures.
$ws.Cells.Item(145, 1).Value = "Somalia"
$ws.Cells.Item(145, 2).Value = 60
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 2
$ws.Cells.Item(145, 5).Value = 56
$ws.Cells.Item(145, 6).Value = 2
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 2

# Row 146: Liberia - shifted down, carries the old row-145 figures.
$ws.Cells.Item(146, 1).Value = "Liberia"
$ws.Cells.Item(146, 2).Value = 59
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 4
$ws.Cells.Item(146, 5).Value = 49
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 6
